$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.179.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6024"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.60%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06945"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2754"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.833.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.737"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009819"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.825.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.567"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -10.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.57%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.877"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.925"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1291"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06512"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.416"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("E31").Value = "  -4.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.779"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.093"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.724"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6459"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.533"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.745"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01754"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.455"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.144.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8876"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.988.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("E46").Value = "  -5.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.612"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.497"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05496"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4537"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.396"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.52%  "
